$d = $word.ActiveDocument

# --- Locate anchors robustly via Find (not hard-coded indices) ---

# Locate "TC 2.1" (the Titolo heading just before the TC 2.1 table).
$rngTitle = $d.Content
$foundTitle = $rngTitle.Find.Execute("TC 2.1", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$titleParaStart = $rngTitle.Paragraphs.Item(1).Range.Start
$titleParaEnd = $rngTitle.Paragraphs.Item(1).Range.End

# The empty "Titolo"-styled paragraph immediately precedes the "TC 2.1" paragraph.
$emptyTitoloStart = $titleParaStart - 1

# Locate the TC 2.1 table (last table in the document).
$tbl = $d.Tables.Item($d.Tables.Count)
$tblStart = $tbl.Range.Start
$tblEnd = $tbl.Range.End

# Locate "RIFERIMENTI" heading paragraph.
$rngRif = $d.Content
$foundRif = $rngRif.Find.Execute("RIFERIMENTI", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rifTextStart = $rngRif.Start
$rifParaEnd = $rngRif.Paragraphs.Item(1).Range.End

# Locate the trailing "Lunghezza massima password..." paragraph.
$rngPwd = $d.Content
$foundPwd = $rngPwd.Find.Execute("Lunghezza massima password", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$pwdParaStart = $rngPwd.Paragraphs.Item(1).Range.Start
$pwdParaEnd = $rngPwd.Paragraphs.Item(1).Range.End
$pwdTextEnd = $pwdParaEnd - 1

# --- Apply deletions from the end of the document backwards so earlier
#     offsets computed above stay valid. ---

# 1) Clear the text runs of the password paragraph, keep the (now empty)
#    paragraph mark / formatting in place.
if ($pwdTextEnd -gt $pwdParaStart) {
    $d.Range($pwdParaStart, $pwdTextEnd).Delete()
}

# 2) Delete the whole "RIFERIMENTI" paragraph (heading text + its mark).
$d.Range($rifTextStart, $rifParaEnd).Delete()

# 3) Delete the run of empty formatted paragraphs between the table and
#    the (now-removed) "RIFERIMENTI" paragraph.
if ($rifTextStart -gt $tblEnd) {
    $d.Range($tblEnd, $rifTextStart).Delete()
}

# 4) Delete the TC 2.1 table itself.
$tbl.Delete()

# 5) Delete the empty "Titolo" paragraph and the "TC 2.1" Titolo paragraph.
$d.Range($emptyTitoloStart, $titleParaEnd).Delete()
